$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Formatting: bold font, thin border all around, centered horizontally,
# top-aligned vertically. Apply fully to B1 first, then replicate the
# exact same style onto A2 by copying formats (keeps a single shared
# cell style in the workbook instead of generating extra intermediate
# style records).
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2

$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)
